$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated weighted variety score (A) and cost/overhead (B) values for BPlocation evaluation.
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 0.2056551268874796

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 0.2056551268874796

$ws.Range("A5").Value = 0.7895487297993826
$ws.Range("B5").Value = 0.2897172436556258

$ws.Range("A6").Value = 0.7581877628298113
$ws.Range("B6").Value = 0.412184793896735

$ws.Range("A7").Value = 0.7581877628298113
$ws.Range("B7").Value = 0.412184793896735

$ws.Range("A8").Value = 0.4717615158673499
$ws.Range("B8").Value = 0.412184793896735

$ws.Range("A9").Value = 0.4717615158673499
$ws.Range("B9").Value = 0.4953213847735023

$ws.Range("A10").Value = 0.2747328878457984
$ws.Range("B10").Value = 0.6028275634437397

$ws.Range("A11").Value = 0.2747328878457984
$ws.Range("B11").Value = 0.6028275634437397

$ws.Range("A12").Value = 0.2747328878457984
$ws.Range("B12").Value = 0.6028275634437397

$ws.Range("A13").Value = 0.2747328878457984
$ws.Range("B13").Value = 0.6028275634437397

$ws.Range("A14").Value = 0.2747328878457984
$ws.Range("B14").Value = 0.6588689746225038

$ws.Range("A15").Value = 0.2747328878457984
$ws.Range("B15").Value = 0.777583435528494

$ws.Range("A16").Value = 0.2747328878457984
$ws.Range("B16").Value = 0.777583435528494

$ws.Range("A17").Value = 0.2747328878457984
$ws.Range("B17").Value = 0.777583435528494

$ws.Range("A18").Value = 0.1802214320095873
$ws.Range("B18").Value = 0.777583435528494

$ws.Range("A19").Value = 0.1454277624071335
$ws.Range("B19").Value = 0.777583435528494

$ws.Range("A20").Value = 0.1454277624071335
$ws.Range("B20").Value = 0.816812423353629

$ws.Range("A21").Value = 0.06537097058055127
$ws.Range("B21").Value = 0.8887917177642474

$ws.Range("A22").Value = 0.06537097058055127
$ws.Range("B22").Value = 0.8887917177642474

$ws.Range("A23").Value = 0.06537097058055127
$ws.Range("B23").Value = 0.8887917177642474

$ws.Range("A24").Value = 0.06272193393914718
$ws.Range("B24").Value = 0.8887917177642474

$ws.Range("A25").Value = 0.04006855027803726
$ws.Range("B25").Value = 0.8887917177642474

$ws.Range("A26").Value = 0.03169260266964554
$ws.Range("B26").Value = 0.8887917177642474

$ws.Range("A27").Value = 0.01644936305047768
$ws.Range("B27").Value = 0.958868974622504
